$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 113, shifting all rows
# from 113..208 down to 114..209 (dimension grows from A1:R208 to A1:R209).
$ws.Rows("113").Insert()

# Populate the newly inserted row 113 with the new record.
$ws.Range("A113").Value = 6
$ws.Range("B113").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C113").Value = "Metropolitana"
$ws.Range("D113").Value = 44658
$ws.Range("E113").Value = 13
$ws.Range("F113").Value = 100112022
$ws.Range("G113").Value = "Arveja Verde"
$ws.Range("H113").Value = "Sin especificar"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 22000
$ws.Range("L113").Value = 23000
$ws.Range("M113").Value = 22575
$ws.Range("N113").Value = "$/saco 25 kilos"
$ws.Range("O113").Value = "Carahue"
$ws.Range("P113").Value = 903
$ws.Range("Q113").Value = 25
$ws.Range("R113").Value = "Hortaliza"
